$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) cell values -------------------------------------------
# Each entry: (1-based column index, header text)
$headers = @(
    @(2, "street"),
    @(3, "city"),
    @(4, "state"),
    @(5, "zip"),
    @(6, "SPATIAL_GEOID"),
    @(7, "social_vulnerability_index"),
    @(8, "gini_inequality_coefficient"),
    @(9, "old_age_dependency_ratio"),
    @(10, "child_dependency_ratio"),
    @(11, "housing_median_year_built"),
    @(12, "housing_percent_occupied_units_lacking_plumbing"),
    @(13, "housing_percent_occupied_lacking_complete_kitchen"),
    @(14, "housing_percent_occupied_units_with_no_bedroom"),
    @(15, "housing_percent_occupied_units_with_no_vehicle_available"),
    @(16, "housing_percent_occupied_units_with_no_computer_included_smartphone"),
    @(17, "housing_percent_occupied_units_with_no_internet_subscription"),
    @(18, "population_density"),
    @(19, "percent_hispanic"),
    @(20, "percent_non_hispanic"),
    @(21, "percent_american_indian_or_alaska_native"),
    @(22, "percent_asian"),
    @(23, "percent_black"),
    @(24, "percent_native_hawaiian_or_other_pacific_islander"),
    @(25, "percent_multiple_race"),
    @(26, "percent_white"),
    @(27, "percent_some_other_race"),
    @(28, "percent_below_100_of_fed_poverty_level"),
    @(29, "percent_households_that_receive_snap"),
    @(30, "percent_households_with_limited_english"),
    @(31, "percent_bachelors_degree_age_25_or_over"),
    @(32, "median_household_income"),
    @(33, "unemployment_rate_age_16_or_over"),
    @(34, "air_quality_indicator_ozone_o3"),
    @(35, "air_quality_indicator_pm25"),
    @(36, "drinking_water_quality_indicator"),
    @(37, "air_quality_indicator_asthma_er_visits"),
    @(38, "food_fraction_of_population_with_low_access"),
    @(39, "food_low_access_tract")
)

foreach ($item in $headers) {
    $col = $item[0]
    $name = $item[1]
    $ws.Cells.Item(1, $col).Value = $name
}

# Copy the existing header style (bold, centered, thin border -- the style
# already used by B1:D1) onto the newly populated header cells so they match
# the pre-existing header formatting instead of getting default formatting.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("E1:AM1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 2 (data) cell values ----------------------------------------------
# Each entry: (1-based column index, value, forceText flag, isNumeric flag)
# forceText -> the value looks numeric (e.g. "24.0") but must be stored as
#              text, so the cell is pre-formatted as Text before the value is
#              assigned (otherwise Excel auto-converts it to a Number and
#              drops formatting such as the trailing ".0").
# isNumeric -> the value is a genuine numeric cell in the target workbook.
$rows2 = @(
    @(2, "1745 T Street Southeast", $false, $false),
    @(3, "Washington", $false, $false),
    @(4, "DC", $false, $false),
    @(5, "20020", $true, $false),
    @(6, "11001007605", $true, $false),
    @(7, "Missing", $false, $false),
    @(8, "0.5317", $true, $false),
    @(9, "24.0", $true, $false),
    @(10, "28.1", $true, $false),
    @(11, "1957", $true, $false),
    @(12, 0, $false, $true),
    @(13, 1.5, $false, $true),
    @(14, "5.1", $true, $false),
    @(15, "45.4", $true, $false),
    @(16, "14.7", $true, $false),
    @(17, "28.1", $true, $false),
    @(19, "6.9", $true, $false),
    @(20, "93.1", $true, $false),
    @(21, "0.0", $true, $false),
    @(22, "1.3", $true, $false),
    @(23, "84.7", $true, $false),
    @(24, "0.0", $true, $false),
    @(25, "2.3", $true, $false),
    @(26, "5.6", $true, $false),
    @(27, "6.0", $true, $false),
    @(28, "28.8", $true, $false),
    @(29, "36.6", $true, $false),
    @(30, "34.2", $true, $false),
    @(31, "292", $true, $false),
    @(32, "40239", $true, $false),
    @(33, "10.3", $true, $false),
    @(34, "Missing", $false, $false),
    @(35, "Missing", $false, $false),
    @(36, "Missing", $false, $false),
    @(37, "Missing", $false, $false),
    @(38, "N/A", $false, $false),
    @(39, "Missing", $false, $false)
)

foreach ($item in $rows2) {
    $col = $item[0]
    $val = $item[1]
    $forceText = $item[2]
    $isNumeric = $item[3]
    if ($forceText) {
        $ws.Cells.Item(2, $col).NumberFormat = "@"
    }
    if ($isNumeric) {
        $ws.Cells.Item(2, $col).Value = [double]$val
    } else {
        $ws.Cells.Item(2, $col).Value = [string]$val
    }
}

# R2 (population_density) has no value for this row in the source data, but
# the column is still present, so write an explicit empty text cell so it
# participates in the used range / dimension like the source workbook does.
# (Assigning "" via .Value clears the cell entirely instead of leaving an
# empty Text-typed cell, so an empty-string formula is used to force Text
# typing with empty content, matching the source workbook's empty inline
# string cell.)
$ws.Cells.Item(2, 18).NumberFormat = "@"
$ws.Cells.Item(2, 18).Formula = "="""""

Write-Output "edit complete"
